$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5750
$ws.Range("J17").Value = 5750
$ws.Range("L17").Value = 17250
$ws.Range("N17").Value = -17586
$ws.Range("H40").Value = 3809.2964
$ws.Range("I40").Value = 2777.077
$ws.Range("J40").Value = 4767.7856
$ws.Range("K40").Value = 2777.077
$ws.Range("L40").Value = 4767.7856
$ws.Range("M40").Value = -2602.077
$ws.Range("N40").Value = -5117.7856
$ws.Range("H103").Value = 2495
$ws.Range("I103").Value = 2990
$ws.Range("K103").Value = 8970
$ws.Range("M103").Value = -8384
$ws.Range("H106").Value = 1900.7273
$ws.Range("I106").Value = 1840.8
$ws.Range("K106").Value = 1840.8
$ws.Range("M106").Value = -1209.8
$ws.Range("H132").Value = 35559.44
$ws.Range("I132").Value = 16486.348
$ws.Range("K132").Value = 49459.04400000001
$ws.Range("M132").Value = -46929.04400000001
$ws.Range("H137").Value = 17243808
$ws.Range("I137").Value = 20002384
$ws.Range("J137").Value = 2715.25
$ws.Range("K137").Value = 60007152
$ws.Range("L137").Value = 8145.75
$ws.Range("M137").Value = -60004602
$ws.Range("N137").Value = -13245.75
$ws.Range("H138").Value = 2238.9624
$ws.Range("I138").Value = 566.7755
$ws.Range("K138").Value = 1700.3265
$ws.Range("M138").Value = 3439.6735

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8507.674000000001
$ws.Range("I32").Value = 8936.977000000001
$ws.Range("K32").Value = 8936.977000000001
$ws.Range("M32").Value = -8649.977000000001
$ws.Range("H45").Value = 2306.375
$ws.Range("I45").Value = 2133.85
$ws.Range("J45").Value = 3169
$ws.Range("K45").Value = 2133.85
$ws.Range("L45").Value = 3169
$ws.Range("M45").Value = -1756.85
$ws.Range("N45").Value = -3923
$ws.Range("H74").Value = 1308.7428
$ws.Range("I74").Value = 1096.1666
$ws.Range("J74").Value = 1533.8235
$ws.Range("K74").Value = 1096.1666
$ws.Range("L74").Value = 1533.8235
$ws.Range("M74").Value = -222.1666
$ws.Range("N74").Value = -3281.8235
$ws.Range("H76").Value = 53798
$ws.Range("J76").Value = 53798
$ws.Range("L76").Value = 53798
$ws.Range("N76").Value = -54474
$ws.Range("H77").Value = 1308.7428
$ws.Range("I77").Value = 1096.1666
$ws.Range("J77").Value = 1533.8235
$ws.Range("K77").Value = 5480.833000000001
$ws.Range("L77").Value = 7669.1175
$ws.Range("M77").Value = -1112.833000000001
$ws.Range("N77").Value = -16405.1175
$ws.Range("H79").Value = 53798
$ws.Range("J79").Value = 53798
$ws.Range("L79").Value = 53798
$ws.Range("N79").Value = -56138
$ws.Range("H132").Value = 1827.5217
$ws.Range("I132").Value = 1849.4
$ws.Range("K132").Value = 5548.200000000001
$ws.Range("M132").Value = -3018.200000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2880.9534
$ws.Range("I134").Value = 2573.4482
$ws.Range("J134").Value = 3517.9285
$ws.Range("K134").Value = 7720.344599999999
$ws.Range("L134").Value = 10553.7855
$ws.Range("M134").Value = -5185.344599999999
$ws.Range("N134").Value = -15623.7855

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2050.6333
$ws.Range("I31").Value = 1885.3462
$ws.Range("K31").Value = 1885.3462
$ws.Range("M31").Value = -1590.3462
$ws.Range("H34").Value = 2050.6333
$ws.Range("I34").Value = 1885.3462
$ws.Range("K34").Value = 1885.3462
$ws.Range("M34").Value = -1683.3462
$ws.Range("H94").Value = 8483.786
$ws.Range("J94").Value = 1649.4546
$ws.Range("L94").Value = 1649.4546
$ws.Range("N94").Value = -2551.4546
$ws.Range("H107").Value = 12316.556
$ws.Range("I107").Value = 1019.8
$ws.Range("K107").Value = 1019.8
$ws.Range("M107").Value = 900.2
$ws.Range("H132").Value = 3553.6956
$ws.Range("I132").Value = 3574.318
$ws.Range("J132").Value = 3100
$ws.Range("K132").Value = 10722.954
$ws.Range("L132").Value = 9300
$ws.Range("M132").Value = -8192.954000000002
$ws.Range("N132").Value = -14360
$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 18630184
$ws.Range("I4").Value = 28743468
$ws.Range("J4").Value = 447.1579
$ws.Range("K4").Value = 86230404
$ws.Range("L4").Value = 1341.4737
$ws.Range("M4").Value = -86230292
$ws.Range("N4").Value = -1565.4737
$ws.Range("H17").Value = 76.666664
$ws.Range("J17").Value = 70
$ws.Range("L17").Value = 210
$ws.Range("N17").Value = -548
$ws.Range("H38").Value = 777.19354
$ws.Range("I38").Value = 79.15000000000001
$ws.Range("K38").Value = 237.45
$ws.Range("M38").Value = 109.55
$ws.Range("H97").Value = 466.66666
$ws.Range("I97").Value = 400
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 1200
$ws.Range("L97").Value = 1500
$ws.Range("N97").Value = -2492
$ws.Range("M97").Value = -704

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 66495.10000000001
$ws.Range("I70").Value = 116271.1
$ws.Range("K70").Value = 116271.1
$ws.Range("M70").Value = -116001.1
$ws.Range("H73").Value = 66495.10000000001
$ws.Range("I73").Value = 116271.1
$ws.Range("K73").Value = 116271.1
$ws.Range("M73").Value = -115335.1
$ws.Range("H102").Value = 2263.2693
$ws.Range("I102").Value = 1422.35
$ws.Range("J102").Value = 5066.3335
$ws.Range("K102").Value = 1422.35
$ws.Range("L102").Value = 5066.3335
$ws.Range("M102").Value = 199.6500000000001
$ws.Range("N102").Value = -8310.333500000001
$ws.Range("H122").Value = 4572.7334
$ws.Range("I122").Value = 2287.889
$ws.Range("K122").Value = 6863.667
$ws.Range("M122").Value = -4413.667
$ws.Range("H132").Value = 3254.1738
$ws.Range("I132").Value = 3461.8
$ws.Range("J132").Value = 2864.875
$ws.Range("K132").Value = 10385.4
$ws.Range("L132").Value = 8594.625
$ws.Range("M132").Value = -7855.400000000001
$ws.Range("N132").Value = -13654.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 715
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H40").Value = 17324.285
$ws.Range("I40").Value = 19170
$ws.Range("J40").Value = 6250
$ws.Range("K40").Value = 19170
$ws.Range("L40").Value = 6250
$ws.Range("M40").Value = -19034
$ws.Range("N40").Value = -6522
$ws.Range("H55").Value = 2168.0908
$ws.Range("I55").Value = 1964.5883
$ws.Range("J55").Value = 2860
$ws.Range("K55").Value = 1964.5883
$ws.Range("L55").Value = 2860
$ws.Range("M55").Value = -1791.5883
$ws.Range("N55").Value = -3206
$ws.Range("H122").Value = 6096.684
$ws.Range("I122").Value = 1835.5
$ws.Range("J122").Value = 7233
$ws.Range("K122").Value = 5506.5
$ws.Range("L122").Value = 21699
$ws.Range("M122").Value = -3056.5
$ws.Range("N122").Value = -26599
$ws.Range("H132").Value = 3494.0925
$ws.Range("I132").Value = 2509.3555
$ws.Range("K132").Value = 7528.066500000001
$ws.Range("M132").Value = -4998.066500000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 49850
$ws.Range("J64").Value = 50133.332
$ws.Range("L64").Value = 50133.332
$ws.Range("N64").Value = -50629.332
$ws.Range("H67").Value = 49850
$ws.Range("J67").Value = 50133.332
$ws.Range("L67").Value = 50133.332
$ws.Range("N67").Value = -51849.332
$ws.Range("H113").Value = 583.6667
$ws.Range("I113").Value = 583.6667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1751.0001
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 418.9999
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 452507.75
$ws.Range("I122").Value = 3141.4
$ws.Range("K122").Value = 9424.200000000001
$ws.Range("M122").Value = -6974.200000000001
$ws.Range("H132").Value = 1773.683
$ws.Range("I132").Value = 1743.05
$ws.Range("K132").Value = 5229.15
$ws.Range("M132").Value = -2699.15
